$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Columns.Item(3).Insert()
$c = $ws.Columns.Item(3)
Write-Output ($c | Get-Member)
